$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I1").Value = "LOA (M)"
$ws.Range("J1").Value = "BREADTH (M)"
$ws.Range("K1").Value = "DEPTH (M)"
$ws.Range("L1").Value = "DRAFT MAX (M)"
$ws.Range("S1").Value = "GT (TON)"
$ws.Range("T1").Value = "KECEPATAN (KNOT)"
$ws.Range("U1").Value = "BOLLARD PULL (TON)"
